# Auto-generated edit script: update cached Leve profit figures per scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 370.5
$ws.Cells.Item(5, 10).Value = 1002
$ws.Cells.Item(5, 12).Value = 1002
$ws.Cells.Item(5, 14).Value = -1232
$ws.Cells.Item(18, 8).Value = 1456
$ws.Cells.Item(18, 9).Value = 1456
$ws.Cells.Item(18, 11).Value = 1456
$ws.Cells.Item(18, 13).Value = -1172
$ws.Cells.Item(40, 8).Value = 1939.091
$ws.Cells.Item(40, 10).Value = 2287.375
$ws.Cells.Item(40, 12).Value = 2287.375
$ws.Cells.Item(40, 14).Value = -2637.375
$ws.Cells.Item(62, 8).Value = 2874.5
$ws.Cells.Item(62, 9).Value = 2879.6
$ws.Cells.Item(62, 11).Value = 2879.6
$ws.Cells.Item(62, 13).Value = -2255.6
$ws.Cells.Item(64, 8).Value = 4422.875
$ws.Cells.Item(64, 9).Value = 3345.75
$ws.Cells.Item(64, 11).Value = 3345.75
$ws.Cells.Item(64, 13).Value = -3097.75
$ws.Cells.Item(65, 8).Value = 2874.5
$ws.Cells.Item(65, 9).Value = 2879.6
$ws.Cells.Item(65, 11).Value = 14398
$ws.Cells.Item(65, 13).Value = -11278
$ws.Cells.Item(67, 8).Value = 4422.875
$ws.Cells.Item(67, 9).Value = 3345.75
$ws.Cells.Item(67, 11).Value = 3345.75
$ws.Cells.Item(67, 13).Value = -2487.75
$ws.Cells.Item(100, 8).Value = 1223.1111
$ws.Cells.Item(100, 9).Value = 401.66666
$ws.Cells.Item(100, 10).Value = 2866
$ws.Cells.Item(100, 11).Value = 401.66666
$ws.Cells.Item(100, 12).Value = 2866
$ws.Cells.Item(100, 13).Value = 139.33334
$ws.Cells.Item(100, 14).Value = -3948
$ws.Cells.Item(106, 8).Value = 1181.25
$ws.Cells.Item(106, 9).Value = 906.3333
$ws.Cells.Item(106, 11).Value = 906.3333
$ws.Cells.Item(106, 13).Value = -275.3333
$ws.Cells.Item(111, 8).Value = 3426.5715
$ws.Cells.Item(111, 9).Value = 1752.75
$ws.Cells.Item(111, 10).Value = 5658.3335
$ws.Cells.Item(111, 11).Value = 5258.25
$ws.Cells.Item(111, 12).Value = 16975.0005
$ws.Cells.Item(111, 13).Value = -2191.25
$ws.Cells.Item(111, 14).Value = -23109.0005
$ws.Cells.Item(116, 8).Value = 3292.7666
$ws.Cells.Item(116, 9).Value = 2629.2942
$ws.Cells.Item(116, 11).Value = 2629.2942
$ws.Cells.Item(116, 13).Value = 812.7058000000002
$ws.Cells.Item(132, 8).Value = 2231.5264
$ws.Cells.Item(132, 9).Value = 2199.9412
$ws.Cells.Item(132, 11).Value = 6599.823600000001
$ws.Cells.Item(132, 13).Value = -4069.823600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 940.3077
$ws.Cells.Item(74, 9).Value = 860.6
$ws.Cells.Item(74, 10).Value = 1206
$ws.Cells.Item(74, 11).Value = 860.6
$ws.Cells.Item(74, 12).Value = 1206
$ws.Cells.Item(74, 13).Value = 13.39999999999998
$ws.Cells.Item(74, 14).Value = -2954
$ws.Cells.Item(77, 8).Value = 940.3077
$ws.Cells.Item(77, 9).Value = 860.6
$ws.Cells.Item(77, 10).Value = 1206
$ws.Cells.Item(77, 11).Value = 4303
$ws.Cells.Item(77, 12).Value = 6030
$ws.Cells.Item(77, 13).Value = 65
$ws.Cells.Item(77, 14).Value = -14766
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).ClearContents()
$ws.Cells.Item(130, 14).Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 7316.8887
$ws.Cells.Item(20, 9).Value = 6984
$ws.Cells.Item(20, 10).Value = 7982.6665
$ws.Cells.Item(20, 11).Value = 6984
$ws.Cells.Item(20, 12).Value = 7982.6665
$ws.Cells.Item(20, 13).Value = -6737
$ws.Cells.Item(20, 14).Value = -8476.666499999999
$ws.Cells.Item(99, 8).Value = 2725.75
$ws.Cells.Item(99, 10).Value = 2833
$ws.Cells.Item(99, 12).Value = 2833
$ws.Cells.Item(99, 14).Value = -5829
$ws.Cells.Item(105, 8).Value = 2389.8
$ws.Cells.Item(105, 9).Value = 2362.25
$ws.Cells.Item(105, 10).Value = 2500
$ws.Cells.Item(105, 11).Value = 2362.25
$ws.Cells.Item(105, 12).Value = 2500
$ws.Cells.Item(105, 13).Value = -615.25
$ws.Cells.Item(105, 14).Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2135.04
$ws.Cells.Item(31, 9).Value = 1525.4445
$ws.Cells.Item(31, 10).Value = 2477.9375
$ws.Cells.Item(31, 11).Value = 1525.4445
$ws.Cells.Item(31, 12).Value = 2477.9375
$ws.Cells.Item(31, 13).Value = -1230.4445
$ws.Cells.Item(31, 14).Value = -3067.9375
$ws.Cells.Item(34, 8).Value = 2135.04
$ws.Cells.Item(34, 9).Value = 1525.4445
$ws.Cells.Item(34, 10).Value = 2477.9375
$ws.Cells.Item(34, 11).Value = 1525.4445
$ws.Cells.Item(34, 12).Value = 2477.9375
$ws.Cells.Item(34, 13).Value = -1323.4445
$ws.Cells.Item(34, 14).Value = -2881.9375
$ws.Cells.Item(53, 8).Value = 39950
$ws.Cells.Item(53, 10).Value = 39950
$ws.Cells.Item(53, 12).Value = 39950
$ws.Cells.Item(53, 14).Value = -41164
$ws.Cells.Item(81, 8).Value = 78450
$ws.Cells.Item(81, 10).Value = 78450
$ws.Cells.Item(81, 12).Value = 78450
$ws.Cells.Item(81, 14).Value = -80446
$ws.Cells.Item(84, 8).Value = 78450
$ws.Cells.Item(84, 10).Value = 78450
$ws.Cells.Item(84, 12).Value = 235350
$ws.Cells.Item(84, 14).Value = -245334
$ws.Cells.Item(88, 8).Value = 19192
$ws.Cells.Item(88, 10).Value = 19192
$ws.Cells.Item(88, 12).Value = 19192
$ws.Cells.Item(88, 14).Value = -20004
$ws.Cells.Item(91, 8).Value = 19192
$ws.Cells.Item(91, 10).Value = 19192
$ws.Cells.Item(91, 12).Value = 19192
$ws.Cells.Item(91, 14).Value = -22000
$ws.Cells.Item(102, 8).Value = 49249.5
$ws.Cells.Item(102, 10).Value = 49249.5
$ws.Cells.Item(102, 12).Value = 49249.5
$ws.Cells.Item(102, 14).Value = -54117.5
$ws.Cells.Item(104, 8).Value = 49966.332
$ws.Cells.Item(104, 10).Value = 49949.5
$ws.Cells.Item(104, 12).Value = 49949.5
$ws.Cells.Item(104, 14).Value = -55191.5
$ws.Cells.Item(107, 8).Value = 1024
$ws.Cells.Item(107, 9).Value = 973.7143
$ws.Cells.Item(107, 10).Value = 1200
$ws.Cells.Item(107, 11).Value = 973.7143
$ws.Cells.Item(107, 12).Value = 1200
$ws.Cells.Item(107, 13).Value = 946.2857
$ws.Cells.Item(107, 14).Value = -5040
$ws.Cells.Item(108, 8).Value = 62752
$ws.Cells.Item(108, 10).Value = 62752
$ws.Cells.Item(108, 12).Value = 62752
$ws.Cells.Item(108, 14).Value = -70432
$ws.Cells.Item(109, 8).Value = 27000
$ws.Cells.Item(109, 10).Value = 27000
$ws.Cells.Item(109, 12).Value = 27000
$ws.Cells.Item(109, 14).Value = -29080
$ws.Cells.Item(114, 8).Value = 40000
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 40000
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).ClearContents()
$ws.Cells.Item(114, 13).Value = 40000
$ws.Cells.Item(114, 14).Value = -48678
$ws.Cells.Item(115, 8).Value = 30000
$ws.Cells.Item(115, 10).Value = 30000
$ws.Cells.Item(115, 12).Value = 30000
$ws.Cells.Item(115, 14).Value = -32350
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).ClearContents()
$ws.Cells.Item(118, 14).Value = 0
$ws.Cells.Item(119, 8).Value = 32500
$ws.Cells.Item(119, 10).Value = 32500
$ws.Cells.Item(119, 12).Value = 32500
$ws.Cells.Item(119, 14).Value = -42176
$ws.Cells.Item(132, 8).Value = 4940
$ws.Cells.Item(132, 9).Value = 5433
$ws.Cells.Item(132, 10).Value = 996
$ws.Cells.Item(132, 11).Value = 16299
$ws.Cells.Item(132, 12).Value = 2988
$ws.Cells.Item(132, 13).Value = -13769
$ws.Cells.Item(132, 14).Value = -8048

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 18999.6
$ws.Cells.Item(68, 8).Value = 1854.32
$ws.Cells.Item(68, 10).Value = 1869.125
$ws.Cells.Item(68, 12).Value = 5607.375
$ws.Cells.Item(68, 14).Value = -7229.375
$ws.Cells.Item(71, 8).Value = 1854.32
$ws.Cells.Item(71, 10).Value = 1869.125
$ws.Cells.Item(71, 12).Value = 16822.125
$ws.Cells.Item(71, 14).Value = -24934.125
$ws.Cells.Item(76, 8).Value = 13001.444
$ws.Cells.Item(76, 9).Value = 6503.25
$ws.Cells.Item(76, 11).Value = 19509.75
$ws.Cells.Item(76, 13).Value = -19126.75
$ws.Cells.Item(79, 8).Value = 13001.444
$ws.Cells.Item(79, 9).Value = 6503.25
$ws.Cells.Item(79, 11).Value = 19509.75
$ws.Cells.Item(79, 13).Value = -18183.75
$ws.Cells.Item(87, 8).Value = 10007
$ws.Cells.Item(87, 9).Value = 10007
$ws.Cells.Item(87, 11).Value = 30021
$ws.Cells.Item(87, 13).Value = -28773
$ws.Cells.Item(90, 8).Value = 10007
$ws.Cells.Item(90, 9).Value = 10007
$ws.Cells.Item(90, 11).Value = 90063
$ws.Cells.Item(90, 13).Value = -83823
$ws.Cells.Item(132, 8).Value = 1999.6666
$ws.Cells.Item(132, 10).Value = 1999
$ws.Cells.Item(132, 12).Value = 17991
$ws.Cells.Item(132, 14).Value = -23051

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2020.375
$ws.Cells.Item(102, 9).Value = 2023.2858
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 11).Value = 2023.2858
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = -401.2858000000001
$ws.Cells.Item(102, 14).Value = -5244
$ws.Cells.Item(132, 8).Value = 3375.8667
$ws.Cells.Item(132, 9).Value = 3438.4285
$ws.Cells.Item(132, 11).Value = 10315.2855
$ws.Cells.Item(132, 13).Value = -7785.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7217
$ws.Cells.Item(7, 9).Value = 4799.4
$ws.Cells.Item(7, 11).Value = 4799.4
$ws.Cells.Item(7, 13).Value = -4687.4
$ws.Cells.Item(9, 8).Value = 2080
$ws.Cells.Item(9, 9).Value = 610
$ws.Cells.Item(9, 11).Value = 610
$ws.Cells.Item(9, 13).Value = -386
$ws.Cells.Item(31, 8).Value = 1395
$ws.Cells.Item(31, 10).Value = 390
$ws.Cells.Item(31, 12).Value = 390
$ws.Cells.Item(31, 14).Value = -886
$ws.Cells.Item(46, 8).Value = 2027.8462
$ws.Cells.Item(46, 9).Value = 1852.6
$ws.Cells.Item(46, 10).Value = 2137.375
$ws.Cells.Item(46, 11).Value = 1852.6
$ws.Cells.Item(46, 12).Value = 2137.375
$ws.Cells.Item(46, 13).Value = -1664.6
$ws.Cells.Item(46, 14).Value = -2513.375
$ws.Cells.Item(55, 8).Value = 1100.6
$ws.Cells.Item(55, 9).Value = 800.6923
$ws.Cells.Item(55, 11).Value = 800.6923
$ws.Cells.Item(55, 13).Value = -627.6923
$ws.Cells.Item(126, 8).Value = 7217
$ws.Cells.Item(126, 9).Value = 4799.4
$ws.Cells.Item(126, 11).Value = 14398.2
$ws.Cells.Item(126, 13).Value = -11928.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1545.75
$ws.Cells.Item(107, 9).Value = 1671.125
$ws.Cells.Item(107, 11).Value = 5013.375
$ws.Cells.Item(107, 13).Value = -3093.375

Write-Output "Applied cached Leve profit updates across all sheets."
